$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Specification")

$ws.Range("B21").Value = "Name and contact information if an agent is being used."
$ws.Range("B25").Value = "Name and contact information if an agent is being used."
$ws.Range("B33").Value = "Telephone number and email address of the applicant."
$ws.Range("B37").Value = "Name and contact information for the parties making the application."
$ws.Range("B43").Value = "Checking whether all the requirements of the form have been met, such as proof of payment or supporting documentation."
$ws.Range("B44").Value = "Details of any conflict of interest that may exist between the applicant and planning authority."
$ws.Range("B47").Value = "Signed and dated verification of the application's accuracy."
$ws.Range("B50").Value = "Further details of any issues relating to trees on the site"
$ws.Range("B53").Value = "Where trees affected by the proposed development are located."
$ws.Range("B63").Value = "Who owns any trees affected by the proposed development."
$ws.Range("B72").Value = "Details of trees affected by the proposed development and what work is being done to them."
